$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Welcome"
$ws.Range("A7").Value = "Today"

$names = @(
  "Vasya",
  "Андрей",
  "Мария",
  "Дрон",
  "Roman",
  "Obama",
  "Алексей",
  "Владимир",
  "Кирилл",
  "Далгат",
  "ﾐ侑ｽﾐｽﾐｾﾐｺﾐｵﾐｽﾑひｸ",
  "Ivanov Ivan",
  "Irtuganov Nickolay",
  "Vasiliev Dmitrie",
  "Qwerty123",
  "Qwerty1234",
  "Qwerty",
  "Qwerty9876",
  "ﾐ籍ｱﾐｰﾐｽﾐｴﾐｾﾐｽ",
  "Абхазия",
  "Alexandr",
  "Alex"
)

$row = 8
foreach ($name in $names) {
    $ws.Cells.Item($row, 1).Value = $name
    $row++
}

$ws.Range("D9").Select()
